$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.751.81"
$ws.Range("E2").Value = "  -3.21%  "
$ws.Range("D3").Value = "2.611.75"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -2.93%  "
$ws.Range("D9").Value = "2.609.10"
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("E10").Value = "  -6.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("E12").Value = "  -5.04%  "
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.57%  "
$ws.Range("D15").Value = "3.084.19"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D17").Value = "63.669.91"
$ws.Range("E17").Value = "  -3.18%  "
$ws.Range("D18").Value = "2.632.30"
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("E21").Value = "  -5.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "343.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.77%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.77%  "
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("E26").Value = "  -4.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "596.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.23%  "
$ws.Range("E29").Value = "  -2.85%  "
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("E33").Value = "  -4.42%  "
$ws.Range("E34").Value = "  -4.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.402"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.41%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("E41").Value = "  -4.82%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.53"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "157.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("E47").Value = "  -5.02%  "
$ws.Range("E48").Value = "  -4.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.630"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.48%  "
$ws.Range("E50").Value = "  -1.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0247"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.79%  "
